$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.689.81"
$ws.Range("E2").Value = "  +1.24%  "

# Row 3
$ws.Range("D3").Value = "1.806.83"
$ws.Range("E3").Value = "  -0.21%  "

# Row 4
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").Value = "'317.61"
$ws.Range("E5").Value = "  -0.17%  "

# Row 6
$ws.Range("E6").Value = "  -0.06%  "

# Row 7
$ws.Range("D7").Value = "'0.5473"
$ws.Range("E7").Value = "  -4.22%  "

# Row 8
$ws.Range("D8").Value = "'0.3805"
$ws.Range("E8").Value = "  -2.09%  "

# Row 9
$ws.Range("D9").Value = "'0.07524"
$ws.Range("E9").Value = "  -0.93%  "

# Row 10
$ws.Range("D10").Value = "'42.43"
$ws.Range("E10").Value = "  -1.14%  "

# Row 11
$ws.Range("E11").Value = "  -2.07%  "

# Row 12
$ws.Range("E12").Value = "  -0.06%  "

# Row 13
$ws.Range("D13").Value = "'20.75"
$ws.Range("E13").Value = "  -2.05%  "

# Row 14
$ws.Range("D14").Value = "'6.167"
$ws.Range("E14").Value = "  -1.58%  "

# Row 15
$ws.Range("D15").Value = "'7.410"
$ws.Range("E15").Value = "  +1.84%  "

# Row 16
$ws.Range("D16").Value = "1.798.29"
$ws.Range("E16").Value = "  -0.67%  "

# Row 17
$ws.Range("D17").Value = "'90.29"
$ws.Range("E17").Value = "  -1.87%  "

# Row 18
$ws.Range("E18").Value = "  -0.48%  "

# Row 19
$ws.Range("D19").Value = "'0.06485"
$ws.Range("E19").Value = "  +0.03%  "

# Row 21
$ws.Range("D21").Value = "'17.39"
$ws.Range("E21").Value = "  +0.46%  "

# Row 22
$ws.Range("D22").Value = "'5.950"
$ws.Range("E22").Value = "  -0.94%  "

# Row 23
$ws.Range("D23").Value = "28.670.84"
$ws.Range("E23").Value = "  +1.11%  "

# Row 24
$ws.Range("D24").Value = "'11.14"
$ws.Range("E24").Value = "  -1.58%  "

# Row 25
$ws.Range("D25").Value = "'2.094"
$ws.Range("E25").Value = "  -2.21%  "

# Row 26
$ws.Range("D26").Value = "'160.08"
$ws.Range("E26").Value = "  +1.19%  "

# Row 27
$ws.Range("E27").Value = "  -1.57%  "

# Row 28
$ws.Range("B28").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C28").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D28").Value = "2.002.30"
$ws.Range("E28").Value = "  -0.88%  "

# Row 29
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "'2.363"
$ws.Range("E29").Value = "  -4.09%  "

# Row 30
$ws.Range("D30").Value = "'123.46"
$ws.Range("E30").Value = "  -0.56%  "

# Row 31
$ws.Range("D31").Value = "'1.121"
$ws.Range("E31").Value = "  -3.70%  "

# Row 32
$ws.Range("D32").Value = "'0.1058"
$ws.Range("E32").Value = "  -1.15%  "

# Row 33
$ws.Range("D33").Value = "'5.655"
$ws.Range("E33").Value = "  -2.38%  "

# Row 34
$ws.Range("D34").Value = "'3.685"
$ws.Range("E34").Value = "  +1.48%  "

# Row 35
$ws.Range("D35").Value = "'0.06683"
$ws.Range("E35").Value = "  +9.05%  "

# Row 36
$ws.Range("E36").Value = "  +1.31%  "

# Row 37
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.02304"
$ws.Range("E37").Value = "  -0.70%  "

# Row 38
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").Value = "'8.751"
$ws.Range("E38").Value = "  -2.61%  "

# Row 39
$ws.Range("D39").Value = "'5.036"
$ws.Range("E39").Value = "  -0.01%  "

# Row 40
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").Value = "'11.30"
$ws.Range("E40").Value = "  -3.33%  "

# Row 41
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.6258"
$ws.Range("E41").Value = "  -2.23%  "

# Row 42
$ws.Range("D42").Value = "'1.200"
$ws.Range("E42").Value = "  +3.01%  "

# Row 43
$ws.Range("D43").Value = "'1.438"
$ws.Range("E43").Value = "  +4.25%  "

# Row 44
$ws.Range("D44").Value = "'13.31"
$ws.Range("E44").Value = "  -1.40%  "

# Row 45
$ws.Range("D45").Value = "'0.5874"
$ws.Range("E45").Value = "  -2.27%  "

# Row 46
$ws.Range("D46").Value = "'3.696"
$ws.Range("E46").Value = "  -0.17%  "

# Row 47
$ws.Range("D47").Value = "'126.82"
$ws.Range("E47").Value = "  +3.33%  "

# Row 48
$ws.Range("D48").Value = "'1.953"
$ws.Range("E48").Value = "  +0.02%  "

# Row 49
$ws.Range("D49").Value = "'1.159"
$ws.Range("E49").Value = "  +0.96%  "

# Row 50
$ws.Range("E50").Value = "  +0.32%  "

# Row 51
$ws.Range("D51").Value = "'72.42"
$ws.Range("E51").Value = "  -1.02%  "
